$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment
$ws.Range("D2").Value = '69.410.34'
$ws.Range("E2").Value = '  -2.58%  '
$ws.Range("D3").Value = '3.693.06'
$ws.Range("E3").Value = '  -3.37%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("E6").Value = '  -5.40%  '
$ws.Range("D7").Value = '3.691.87'
$ws.Range("E7").Value = '  -3.40%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -4.43%  '
$ws.Range("E10").Value = '  -8.26%  '
$ws.Range("E11").Value = '  -2.14%  '
$ws.Range("E12").Value = '  -3.79%  '
$ws.Range("E13").Value = '  -5.08%  '
$ws.Range("E14").Value = '  -6.91%  '
$ws.Range("D15").Value = '4.310.55'
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("D16").Value = '3.689.83'
$ws.Range("E16").Value = '  -3.09%  '
$ws.Range("D17").Value = '69.454.50'
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").Value = '  -6.93%  '
$ws.Range("E20").Value = '  -7.63%  '
$ws.Range("E21").Value = '  -5.72%  '
$ws.Range("E22").Value = '  -6.53%  '
$ws.Range("E23").Value = '  -7.70%  '
$ws.Range("E24").Value = '  -4.99%  '
$ws.Range("D25").Value = '3.835.52'
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("E26").Value = '  -9.49%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -5.04%  '
$ws.Range("E29").Value = '  -7.90%  '
$ws.Range("E31").Value = '  -9.80%  '
$ws.Range("E32").Value = '  -7.42%  '
$ws.Range("E33").Value = '  -7.18%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E35").Value = '  -7.07%  '
$ws.Range("E36").Value = '  -4.21%  '
$ws.Range("D37").Value = '3.657.63'
$ws.Range("E37").Value = '  -3.33%  '
$ws.Range("E38").Value = '  -6.90%  '
$ws.Range("E39").Value = '  +5.47%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("E41").Value = '  -7.70%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -7.08%  '
$ws.Range("E45").Value = '  -4.61%  '
$ws.Range("E46").Value = '  -3.27%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E48").Value = '  -14.60%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("E49").Value = '  -7.26%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("E51").Value = '  +0.20%  '

# Numeric-looking price strings: must stay as TEXT (matches source data which
# stores these as inline strings, not numbers), so force Text format first.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '692.74'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.00'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.56'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.34'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '482.65'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.03'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.668'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.06'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.46'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.74'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.88'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.12'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.166'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.55'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.35'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.33'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0935'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.953'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '164.01'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '48.04'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.16'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000288'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
